$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# Update the disclaimer date text in A41 (shared string with CONFIDENTIAL notice)
$ws.Range("A41").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-18 for illustrative purposes only and are subject to change."
$ws.Rows(41).AutoFit()

# Update Weight (D) and Percent Change (E) values for rows 2-37
$ws.Range("D2").Value = 0.02885720993856917
$ws.Range("E2").Value = -0.01209063214013706
$ws.Range("D3").Value = 0.0285932617341249
$ws.Range("E3").Value = -0.02623001707279216
$ws.Range("D4").Value = 0.03143160307430976
$ws.Range("E4").Value = -0.05458066077695745
$ws.Range("D5").Value = 0.06626663755980836
$ws.Range("E5").Value = -0.03435882553663749
$ws.Range("D6").Value = 0.01400573839034706
$ws.Range("E6").Value = -0.001131648434553045
$ws.Range("D7").Value = 0.01415049781392137
$ws.Range("E7").Value = -0.03113799283154106
$ws.Range("D8").Value = 0.03097259947138508
$ws.Range("E8").Value = 0.0006140745895932831
$ws.Range("D9").Value = 0.03215983807232012
$ws.Range("E9").Value = -0.04021553423577362
$ws.Range("D10").Value = 0.03324859799969949
$ws.Range("E10").Value = -0.03746853786896498
$ws.Range("D11").Value = 0.02958671290233345
$ws.Range("E11").Value = 0.01324247878632034
$ws.Range("D12").Value = 0.01631449270063814
$ws.Range("E12").Value = -0.04290155440414511
$ws.Range("D13").Value = 0.01684027583327011
$ws.Range("E13").Value = -0.02933942375263532
$ws.Range("D14").Value = 0.01566444892412049
$ws.Range("E14").Value = -0.01343694350008096
$ws.Range("D15").Value = 0.008375336459995997
$ws.Range("E15").Value = -0.03431570448122723
$ws.Range("D16").Value = 0.007293339016637612
$ws.Range("E16").Value = -0.04357904496986564
$ws.Range("D17").Value = 0.03152670050585494
$ws.Range("E17").Value = -0.01707957957957973
$ws.Range("D18").Value = 0.0310218388081629
$ws.Range("E18").Value = -0.0146122143124765
$ws.Range("D19").Value = 0.03178156162239599
$ws.Range("E19").Value = -0.02084580091761434
$ws.Range("D20").Value = 0.03000957948127098
$ws.Range("E20").Value = -0.01897820499278202
$ws.Range("D21").Value = 0.04419029714566115
$ws.Range("E21").Value = -0.02623524685808287
$ws.Range("D22").Value = 0.02915898578800584
$ws.Range("E22").Value = 0.008841861139295526
$ws.Range("D23").Value = 0.02958903750621567
$ws.Range("E23").Value = 0.01117737385280182
$ws.Range("D24").Value = 0.03041215860814554
$ws.Range("E24").Value = 0.01369606003752333
$ws.Range("D25").Value = 0.01546072909298817
$ws.Range("E25").Value = -0.03324220885729912
$ws.Range("D26").Value = 0.01479293379591541
$ws.Range("E26").Value = -0.03107142857142853
$ws.Range("D27").Value = 0.02999436389222375
$ws.Range("E27").Value = 0.007982639696194838
$ws.Range("D28").Value = 0.03161630341913305
$ws.Range("E28").Value = 0.005240361477995892
$ws.Range("D29").Value = 0.03005586023128963
$ws.Range("E29").Value = -0.02666216672291599
$ws.Range("D30").Value = 0.03060446674749243
$ws.Range("E30").Value = -0.01139345394282554
$ws.Range("D31").Value = 0.03383249622938683
$ws.Range("E31").Value = -0.04637871263937043
$ws.Range("D32").Value = 0.03106135707416057
$ws.Range("E32").Value = -0.008157461457865778
$ws.Range("D33").Value = 0.03011672258747854
$ws.Range("E33").Value = -0.03587066352307178
$ws.Range("D34").Value = 0.02977014528140281
$ws.Range("E34").Value = 0.02802544153557851
$ws.Range("D35").Value = 0.03099669082070985
$ws.Range("E35").Value = -0.0001159017153454789
$ws.Range("D36").Value = 0.02949732131668099
$ws.Range("E36").Value = 0.01650654458701406
$ws.Range("D37").Value = 0.03074986015394372
$ws.Range("E37").Value = -0.0171605684910795
$ws.Range("E38").Value = -0.01589521869133714

$ws.Protect()
